$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $val) {
    $rng = $ws.Range($cellRef)
    $rng.NumberFormat = "@"
    $rng.Value = $val
    $rng.Style = "Normal"
}

Set-TextValue 'D2' '46.942.69'
$ws.Range('E2').Value = '  +6.27%  '

Set-TextValue 'D3' '2.328.49'
$ws.Range('E3').Value = '  +5.11%  '

$ws.Range('E4').Value = '  -0.48%  '

Set-TextValue 'D5' '304.70'
$ws.Range('E5').Value = '  +1.11%  '

Set-TextValue 'D6' '97.42'
$ws.Range('E6').Value = '  +10.07%  '

Set-TextValue 'D7' '0.579'
$ws.Range('E7').Value = '  +4.96%  '

Set-TextValue 'D8' '1.00'
$ws.Range('E8').Value = '  -0.53%  '

Set-TextValue 'D9' '0.539'
$ws.Range('E9').Value = '  +9.91%  '

Set-TextValue 'D10' '35.98'
$ws.Range('E10').Value = '  +7.94%  '

Set-TextValue 'D11' '0.0812'
$ws.Range('E11').Value = '  +4.66%  '

Set-TextValue 'D12' '7.47'
$ws.Range('E12').Value = '  +9.00%  '

$ws.Range('E13').Value = '  +0.95%  '

Set-TextValue 'D14' '2.684.39'
$ws.Range('E14').Value = '  +5.07%  '

Set-TextValue 'D15' '2.335.68'
$ws.Range('E15').Value = '  +1.78%  '

Set-TextValue 'D16' '0.839'
$ws.Range('E16').Value = '  +5.24%  '

Set-TextValue 'D17' '14.14'
$ws.Range('E17').Value = '  +8.13%  '

Set-TextValue 'D18' '46.835.60'
$ws.Range('E18').Value = '  +6.56%  '

Set-TextValue 'D19' '13.59'
$ws.Range('E19').Value = '  +21.54%  '

Set-TextValue 'D20' '0.0₃0956'
$ws.Range('E20').Value = '  +6.29%  '

Set-TextValue 'D21' '6.21'
$ws.Range('E21').Value = '  +4.42%  '

Set-TextValue 'D22' '67.95'
$ws.Range('E22').Value = '  +6.01%  '

Set-TextValue 'D23' '254.79'
$ws.Range('E23').Value = '  +9.97%  '

Set-TextValue 'D24' '2.98'
$ws.Range('E24').Value = '  +4.60%  '

Set-TextValue 'D25' '2.01'
$ws.Range('E25').Value = '  +6.01%  '

$ws.Range('E26').Value = '  -0.21%  '

Set-TextValue 'D27' '42.17'
$ws.Range('E27').Value = '  +17.10%  '

Set-TextValue 'D28' '2.28'
$ws.Range('E28').Value = '  +1.36%  '

Set-TextValue 'D29' '9.95'
$ws.Range('E29').Value = '  +6.20%  '

Set-TextValue 'D30' '20.34'
$ws.Range('E30').Value = '  +5.28%  '

$ws.Range('E31').Value = '  +5.44%  '

Set-TextValue 'D32' '0.0817'
$ws.Range('E32').Value = '  +8.92%  '

Set-TextValue 'D33' '147.92'
$ws.Range('E33').Value = '  +1.40%  '

$ws.Range('E34').Value = '  +2.26%  '

$ws.Range('B35').Value = 'LidoDAOToken'
$ws.Range('C35').Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
Set-TextValue 'D35' '3.13'
$ws.Range('E35').Value = '  +7.38%  '

$ws.Range('B36').Value = 'Kaspa'
$ws.Range('C36').Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
Set-TextValue 'D36' '0.114'
$ws.Range('E36').Value = '  +8.80%  '

$ws.Range('E37').Value = '  +3.96%  '

Set-TextValue 'D38' '1.82'
$ws.Range('E38').Value = '  +6.77%  '

Set-TextValue 'D39' '4.00'
$ws.Range('E39').Value = '  +11.19%  '

$ws.Range('E40').Value = '  +9.38%  '

Set-TextValue 'D41' '3.41'
$ws.Range('E41').Value = '  +7.06%  '

Set-TextValue 'D42' '14.07'
$ws.Range('E42').Value = '  -0.29%  '

Set-TextValue 'D43' '0.999'
$ws.Range('E43').Value = '  -0.70%  '

$ws.Range('E44').Value = '  +18.96%  '

Set-TextValue 'D45' '92.30'
$ws.Range('E45').Value = '  +18.32%  '

Set-TextValue 'D46' '1.802.93'
$ws.Range('E46').Value = '  +3.91%  '

$ws.Range('B47').Value = 'Algorand'
$ws.Range('C47').Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
Set-TextValue 'D47' '0.195'
$ws.Range('E47').Value = '  +9.10%  '

$ws.Range('B48').Value = 'ordi'
$ws.Range('C48').Value = 'https://coinranking.com/coin/j7-7vPrOi+ordi-ordi'
Set-TextValue 'D48' '74.29'
$ws.Range('E48').Value = '  +12.25%  '

Set-TextValue 'D49' '98.89'
$ws.Range('E49').Value = '  +4.95%  '

$ws.Range('B50').Value = 'THORChain'
$ws.Range('C50').Value = 'https://coinranking.com/coin/ybmU-kKU+thorchain-rune'
Set-TextValue 'D50' '4.89'
$ws.Range('E50').Value = '  +5.91%  '

$ws.Range('B51').Value = 'MultiversX'
$ws.Range('C51').Value = 'https://coinranking.com/coin/omwkOTglq+multiversx-egld'
Set-TextValue 'D51' '55.26'
$ws.Range('E51').Value = '  +7.22%  '
